$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append row 45 with the new test-mail log entry ---
$ws.Range("A45").Value = "Kun je dit intern overleggen?"
$ws.Range("B45").Value = "mailmind.test@zohomail.eu"
$ws.Range("C45").Value = "Testmail #13: Kun je dit intern overleggen?"
$ws.Range("D45").Value = "Overig"
$ws.Range("E45").Value = "Beste afzender,`nBedankt voor je e-mail. Kun je wat meer specifieke informatie geven over waarover je precies wilt dat er intern overlegd wordt? Op die manier kan ik ervoor zorgen dat je aanvraag bij de juiste persoon of afdeling terechtkomt.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("F45").Value = "2025-06-26 23:47:16"
$ws.Range("G45").Value = "Ja"
$ws.Range("H45").Value = "Nee"
$ws.Range("I45").Value = "Ja"

# Re-settle the row height back to the sheet default — entering the
# multi-line "Antwoord" text auto-expands the row, but the rest of the
# log rows keep the default (non-custom) height.
$ws.Rows.Item(45).AutoFit()

# --- Extend the conditional-formatting ranges to cover the new row ---
$ranges = @("D2:D44", "G2:G44", "H2:H44", "I2:I44")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newAddr = "$col" + "2:" + "$col" + "45"
    $fcs = $ws.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($newAddr))
    }
}

# --- Dashboard: "Overig" category count goes from 1 to 2 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B7").Value = 2
